# Rust Report.pptx - title slide text/box update
# - Reposition/resize the title textbox (shape 3 on slide 1)
# - Shrink the title font from 54pt to 48pt
# - Change the Chinese portion of the title text

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)

# Resize/reposition the title shape (values in points; 1 pt = 12700 EMU)
# off x: 615315 -> 465455 EMU  (48.45pt -> 36.65pt)
# ext cx: 10645775 -> 11260455 EMU (838.25pt -> 886.65pt)
$shp.Left = 36.65
$shp.Width = 886.65

$tr = $shp.TextFrame.TextRange

# Replace the Chinese run's text (characters 10-18 = "集成工作阶段性汇报")
# while keeping its existing run-level formatting (typeface, lang, bold, ...).
$zhRun = $tr.Characters(10, 9)
$zhRun.Text = "指令集支持工作阶段性汇报"

# Shrink the whole title text from 54pt to 48pt.
$tr.Font.Size = 48
